$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the header style (s="1") to column A, rows 2-8 (A9 will be removed below).
$ws.Range("A1").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

# --- Row 2 (ID = 1) ---
$ws.Range("C2").Value = 3
$ws.Range("F2").Value = "Prévention"
$ws.Range("G2").Value = "Système (Machine)"
$ws.Range("H2").Value = "Bureautique : Word, Excel.. "
$ws.Range("I2").Value = "Assez Importante"
$ws.Range("J2").Value = "test"

# --- Row 3 (ID = 2) ---
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 5
$ws.Range("J3").Value = ""

# --- Row 4 (ID = 3) ---
$ws.Range("C4").Value = 56
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("N4").Value = ""

# --- Row 5 (ID = 4) ---
$ws.Range("C5").Value = 66
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("N5").Value = ""

# --- Row 6 (ID = 5) ---
$ws.Range("C6").Value = 74
$ws.Range("J6").Value = "rapide"

# --- Row 7 (ID = 6) ---
$ws.Range("C7").Value = 85
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 7
$ws.Range("J7").Value = ""

# --- Row 8 (ID = 7) ---
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "Test"
$ws.Range("I8").Value = "Neutre"
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = "Neutre"
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = "Neutre"
$ws.Range("O8").Value = ""

# --- Row 9 removed entirely ---
$ws.Rows(9).Delete()
